$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.721.90"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").Value = "1.592.53"
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "209.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.35%  "
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.33"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.252"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.21%  "
$ws.Range("E10").Value = "  +0.27%  "
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("D12").Value = "1.819.24"
$ws.Range("E12").Value = "  -0.31%  "
$ws.Range("D13").Value = "1.583.33"
$ws.Range("E13").Value = "  -1.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.84"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.529"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.28%  "
$ws.Range("D16").Value = "27.722.36"
$ws.Range("E16").Value = "  +0.56%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.29"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "218.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("D19").Value = "0.0₃0696"
$ws.Range("E19").Value = "  +0.61%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.06%  "
$ws.Range("E21").Value = "  +0.23%  "
$ws.Range("E22").Value = "  -1.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.79"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.38%  "
$ws.Range("E24").Value = "  -2.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.83"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.09"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.90%  "
$ws.Range("E27").Value = "  +0.22%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("E30").Value = "  +0.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0475"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.80%  "
$ws.Range("E32").Value = "  -2.24%  "
$ws.Range("D33").Value = "1.383.00"
$ws.Range("E33").Value = "  +0.93%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.98"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.54"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.55%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.970"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.34"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.85%  "
$ws.Range("E38").Value = "  +2.58%  "
$ws.Range("E39").Value = "  -0.48%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.828"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.29%  "
$ws.Range("E41").Value = "  +0.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.988"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.61%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "64.54"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.94%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.29%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.27"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.94%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.75"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.17%  "
$ws.Range("D47").Value = "1.730.92"
$ws.Range("E47").Value = "  -0.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.05"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.92%  "
$ws.Range("E49").Value = "  +0.60%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0967"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0496"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.39%  "